# Apply "Added p2processanalyse and peter-assesments" edit:
# Fill in Peter's self/peer assessment on the first worksheet
# ("Peer  and self assessment").
#
# NOTE: new text values are written in a specific order below so that the
# workbook's shared-string table is built up in the same order as the
# target file (matching the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 header label becomes "Self assesment (Peter)"
$ws.Range("A2").Value = "Self assesment (Peter)"

# Grades already exist in the shared strings ("Excellent" / "Insufficient"),
# so set those first for each row, then set the new comment text.

$ws.Range("B2").Value = "Excellent"
$ws.Range("B3").Value = "Excellent"
$ws.Range("C3").Value = "Good use of communication programs, and making sure everybody is on the right track. "

$ws.Range("B4").Value = "Excellent"
$ws.Range("C4").Value = "Very motivated. Asking everyone where we are in the project, and what needs to be done. "

$ws.Range("B5").Value = "Excellent"
$ws.Range("C5").Value = "Guiding everyone on how to use the programs (github/sourcetree etc..) very helpful and good overview."

$ws.Range("B6").Value = "Excellent"
$ws.Range("C6").Value = "Same as Ahmed. Asking relevant questions and keeping close contact. "

$ws.Range("B7").Value = "Excellent"
$ws.Range("C7").Value = "Keeping close contact, and bringing a lot of ideas to the table."

$ws.Range("B8").Value = "Excellent"
$ws.Range("C8").Value = "Keeping contact and helps create overall overview of the project."

$ws.Range("B9").Value = "Insufficient"
$ws.Range("C9").Value = "Haven't participated at all."

$ws.Range("B10").Value = "Excellent"
$ws.Range("C10").Value = "Same as Ahmed. Asking relevant questions and keeping close contact. "

$ws.Range("C2").Value = "Keeping close contact (asking/responding when needed). Also very motivated and helps keeping an overview"

# Update the selection shown when the sheet is reopened
$ws.Range("C10").Select()
